$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "VENTAS POR GRUPO"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Column J (10th column) width 9 -> 12 (ColumnWidth property pads by ~0.8333,
# so subtract that offset to land on the exact stored width of 12)
$ws1.Columns.Item(10).ColumnWidth = 11.166666666666666

# New sales figures recorded for this period
$ws1.Range("C4").Value = 1555.2
$ws1.Range("M4").Value = 2004.28
$ws1.Range("H10").Value = 908.1
$ws1.Range("I10").Value = 208.8
$ws1.Range("M35").Value = 4270.55
$ws1.Range("J37").Value = 258.08

# Row 55 "count of advisors with sales > 0 out of 53" summary labels
$ws1.Range("C55").Value = "1 de 53"
$ws1.Range("H55").Value = "1 de 53"
$ws1.Range("I55").Value = "6 de 53"
$ws1.Range("J55").Value = "1 de 53"
$ws1.Range("M55").Value = "6 de 53"

# ---------------------------------------------------------------------------
# Sheet "VENTA MENSUAL"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F4").Value = 3559.48
$ws2.Range("F10").Value = 1116.9
$ws2.Range("F35").Value = 4270.55
$ws2.Range("F37").Value = 258.08
$ws2.Range("F59").Value = 17988.2

# ---------------------------------------------------------------------------
# Sheet "CUMPLIMIENTO MENSUAL"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Column widths: E 22 -> 23, F 25 -> 24
$ws3.Columns.Item(5).ColumnWidth = 22.166666666666668
$ws3.Columns.Item(6).ColumnWidth = 23.166666666666668

$ws3.Range("D2").Value = 1555.2
$ws3.Range("E2").Value = 4642.38402943659
$ws3.Range("F2").Value = 0.2509364927709387

$ws3.Range("D6").Value = 908.1
$ws3.Range("E6").Value = 1999.48368146026
$ws3.Range("F6").Value = 0.3123211915757932

$ws3.Range("D7").Value = 756.9
$ws3.Range("E7").Value = 129.811016287574
$ws3.Range("F7").Value = 0.8536039206650904

$ws3.Range("D9").Value = 230.43
$ws3.Range("E9").Value = -230.43

$ws3.Range("D12").Value = 8074.13
$ws3.Range("E12").Value = 44588.99000000001
$ws3.Range("F12").Value = 0.1533165904336849

$ws3.Range("D14").Value = 15967.51
$ws3.Range("E14").Value = 83048.99661190616
$ws3.Range("F14").Value = 0.1612610921791499
